# heating costs do not matter for housing benefit!
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the heating-cost term (+I<row>) from the AP column formula (rows 2-13)
for ($row = 2; $row -le 13; $row++) {
    $ws.Range("AP$row").Formula = "=ROUND(MAX(MIN(H$row,AM$row),AN$row)+4,-1)-5"
}

# T7 no longer derives from a formula; replace with a plain value
$ws.Range("T7").Value = 16000

# Update the view: scroll so column N is the left-most visible column,
# and move the active selection to S7
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("S7").Select()
